$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E. This shifts the existing D:K quarterly
# columns to F:M, matching the new layout (two additional quarters of
# data are being added at the front of each table, most-recent-first).
$ws.Range("D1:E1").EntireColumn.Insert()

# The freshly inserted D:E columns inherit column C's formatting by
# default; copy the (now-shifted) column F's formats across D:E so the
# new columns match the rest of the data table (number formats, fonts,
# alignment) for every row in the sheet's used range.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

function SetRow([int]$r, [object[]]$vals) {
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("D$r" + ":M$r").Value = $arr
}

SetRow 7 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
SetRow 8 @(178700,185800,218300,190700,181300,272500,180300,167700,180400,174300)
SetRow 9 @(6900,5800,31400,8400,6500,63400,9600,700,9000,2500)
SetRow 10 @(171800,180000,186900,182300,174800,209100,170700,167000,171400,171800)
SetRow 11 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 12 @("NA","NA","NA","NA","NA","NA","NA","NA","NA","NA")
SetRow 13 @(0,0,0,0,0,0,0,0,0,0)
SetRow 14 @(-39300,400,200,0,2100,1000,900,300,100,-6600)
SetRow 15 @(47400,51500,51500,55700,55300,55400,52100,49700,50900,50000)
SetRow 16 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 17 @(130700,153600,201800,162900,187700,214600,168300,126800,127100,127800)
SetRow 18 @(48000,32200,16500,27800,-6400,57900,12000,40900,53300,46500)
SetRow 19 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 20 @(66100,42800,247400,27500,153400,23900,122800,55900,53700,23400)
SetRow 21 @(161500,126500,315400,111000,202300,137200,186900,146500,157900,119900)
SetRow 22 @(56900,55200,67200,58900,58800,75700,104200,100000,49900,51300)
SetRow 23 @(57200,19800,196700,-3600,88200,6100,30600,-3200,57100,18600)
SetRow 24 @(24300,6900,29400,-2600,27600,-3700,8700,-4100,11900,5500)
SetRow 25 @(0,0,0,0,0,0,0,0,0,0)
SetRow 26 @(32900,12900,167300,-1000,60600,9800,21900,900,45200,13100)
SetRow 27 @(30600,11900,108800,-2600,54000,-9100,9200,500,14000,-2900)
SetRow 28 @(0,0,0,0,0,0,0,0,0,0)
SetRow 29 @(0,"NA","NA","NA",44800,"NA","NA","NA","NA","NA")
SetRow 30 @(0,0,0,0,0,0,0,0,0,0)
SetRow 31 @(0,0,0,0,0,0,0,0,0,0)
SetRow 32 @(-66100,-42800,-247400,-27500,-153400,-23900,-122800,-55900,-53700,-23400)
SetRow 33 @(30600,11900,108800,-2600,98800,-9100,9200,500,14000,-2900)
SetRow 34 @(0,0,0,0,0,0,0,0,0,0)
SetRow 35 @(30600,11900,108800,-2600,98800,-9100,9200,500,14000,-2900)
SetRow 38 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
SetRow 39 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 40 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 41 @(488000,419400,447100,433000,351300,512900,569200,168000,260200,282700)
SetRow 42 @(0,0,0,0,0,0,0,0,0,0)
SetRow 43 @(56600,75800,64600,48300,62700,86600,93400,80200,71300,64100)
SetRow 44 @(0,0,0,0,0,0,0,0,0,0)
SetRow 45 @(13700,17800,16800,14500,13300,12300,12900,11800,10200,14200)
SetRow 46 @(0,0,0,0,0,0,0,0,0,0)
SetRow 47 @(887700,819100,642600,573800,578600,1320800,1289700,1298400,1280800,1240000)
SetRow 48 @(5733900,5853900,6133500,6713100,6488000,6344400,6085600,5855600,5839600,6155000)
SetRow 49 @(61100,70100,77900,83700,87000,91200,92600,93300,96300,110800)
SetRow 50 @(0,0,0,0,0,0,0,0,0,0)
SetRow 51 @(0,0,0,0,0,0,0,0,0,0)
SetRow 52 @(0,0,0,0,0,0,0,"NA",28400,32800)
SetRow 53 @(0,0,0,0,0,0,0,0,0,0)
SetRow 54 @(7357100,7381400,7505400,7962800,7724800,8536600,8284300,7626000,7659100,7966800)
SetRow 55 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 56 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 57 @(513100,534300,516900,504800,485400,513400,452200,392500,423300,454700)
SetRow 58 @(0,0,0,0,0,0,0,0,0,0)
SetRow 59 @(0,0,0,0,0,0,0,0,0,0)
SetRow 60 @(0,0,0,0,0,0,0,0,0,0)
SetRow 61 @(5412800,5405300,5538800,5891300,5661900,5628900,5480400,4896300,4892700,4984100)
SetRow 62 @(0,0,0,0,0,0,0,0,0,0)
SetRow 63 @(0,0,0,0,0,0,0,0,0,0)
SetRow 64 @(0,0,0,0,0,0,0,0,0,0)
SetRow 65 @(0,0,0,0,0,0,0,0,0,0)
SetRow 66 @(6110400,6128400,6235100,6617400,6359200,7523200,7256800,6601000,6611100,6903600)
SetRow 67 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 68 @(0,0,0,0,0,0,0,0,0,0)
SetRow 69 @(0,0,0,0,0,0,0,0,0,0)
SetRow 70 @(0,0,0,0,0,0,0,0,0,0)
SetRow 71 @(0,0,0,0,0,0,0,0,0,0)
SetRow 72 @(-56400,-56700,-40700,-121500,-90600,-160400,-132000,-121900,-112200,-106700)
SetRow 73 @(0,0,0,0,0,0,0,0,0,0)
SetRow 74 @(0,0,0,0,0,0,0,0,0,0)
SetRow 75 @(0,0,0,0,0,0,0,0,0,0)
SetRow 76 @(1246700,1253000,1270300,1345400,1365600,1013400,1027500,1025000,1048000,1063200)
SetRow 77 @(0,0,0,0,0,0,0,0,0,0)
SetRow 80 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
SetRow 81 @(30600,11900,108800,-2600,98800,-9100,9200,500,14000,-2900)
SetRow 82 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 83 @(47400,51500,51500,55700,55300,55400,52100,49700,50900,50000)
SetRow 84 @(0,0,0,0,0,0,0,0,0,0)
SetRow 85 @(0,0,0,0,0,0,0,0,0,0)
SetRow 86 @(0,0,0,0,0,0,0,0,0,0)
SetRow 87 @(0,0,0,0,0,0,0,0,0,0)
SetRow 88 @(0,0,0,0,0,0,0,0,0,0)
SetRow 89 @(1300,32100,10900,48800,-56600,61400,57000,11200,24500,67900)
SetRow 90 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 91 @(-6900,-5800,-8000,-8400,-6400,271900,-237900,-47000,-201200,-300900)
SetRow 92 @(0,0,0,0,0,0,0,0,0,0)
SetRow 93 @(0,0,0,0,0,0,0,0,0,0)
SetRow 94 @(80500,84100,483100,-54600,72500,-141800,8300,-9200,36700,-204000)
SetRow 95 @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
SetRow 96 @(-27000,-26700,-28200,-29300,-2100,-19200,-19400,-18500,-17300,-16500)
SetRow 97 @(0,0,0,0,0,0,0,0,0,0)
SetRow 98 @(0,0,0,0,0,0,0,0,0,0)
SetRow 99 @(0,0,0,0,0,0,0,0,0,0)
SetRow 100 @(-8000,-137000,-462900,79100,-873900,28300,368500,-88200,-116100,262900)
SetRow 101 @(-5200,-6900,-17000,8400,-23400,20000,22700,8800,-34600,-7600)
SetRow 102 @(68600,-27700,14100,81700,-881400,-32100,456500,-77400,-89500,119200)
